$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 2).Value = "Optical Density"
}
